# Adds a "LandLeaf" breakdown table (columns AB:AH, rows 8-22) to Sheet1.
# New data: per-LandLeaf-category area change (thousand km^2) next to the
# existing tables, plus a "Land for Food" = Animal Feed + Crops summary row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header (row 8): year / LandLeaf / Units labels for the new block
$ws.Range("AC8").Value = 2050
$ws.Range("AD8").Value = "LandLeaf"
$ws.Range("AE8").Value = "Units"

# Row 9 - Animal Feed
$ws.Range("AB9").Value = 0
$ws.Range("AC9").Value = -42.537666999999999
$ws.Range("AD9").Value = "Animal"
$ws.Range("AE9").Value = "Feed"
$ws.Range("AF9").Value = "thousand"
$ws.Range("AG9").Value = 'km$^2$'

# Row 10 - Biomass for Energy
$ws.Range("AB10").Value = 1
$ws.Range("AC10").Value = 52.946877999999998
$ws.Range("AD10").Value = "Biomass"
$ws.Range("AE10").Value = "for"
$ws.Range("AF10").Value = "Energy"
$ws.Range("AG10").Value = "thousand"
$ws.Range("AH10").Value = 'km$^2$'

# Row 11 - Crops
$ws.Range("AB11").Value = 2
$ws.Range("AC11").Value = -80.947716
$ws.Range("AD11").Value = "Crops"
$ws.Range("AE11").Value = "thousand"
$ws.Range("AF11").Value = 'km$^2$'

# Row 12 - Forest
$ws.Range("AB12").Value = 3
$ws.Range("AC12").Value = 444.701007
$ws.Range("AD12").Value = "Forest"
$ws.Range("AE12").Value = "thousand"
$ws.Range("AF12").Value = 'km$^2$'

# Row 13 - Grasslands
$ws.Range("AB13").Value = 4
$ws.Range("AC13").Value = 92.432309000000004
$ws.Range("AD13").Value = "Grasslands"
$ws.Range("AE13").Value = "thousand"
$ws.Range("AF13").Value = 'km$^2$'

# Row 14 - Other Arable Land
$ws.Range("AB14").Value = 5
$ws.Range("AC14").Value = -123.79096800000001
$ws.Range("AD14").Value = "Other"
$ws.Range("AE14").Value = "Arable"
$ws.Range("AF14").Value = "Land"
$ws.Range("AG14").Value = "thousand"
$ws.Range("AH14").Value = 'km$^2$'

# Row 15 - Pasture
$ws.Range("AB15").Value = 6
$ws.Range("AC15").Value = -322.43952899999999
$ws.Range("AD15").Value = "Pasture"
$ws.Range("AE15").Value = "thousand"
$ws.Range("AF15").Value = 'km$^2$'

# Row 16 - Rock and Desert
$ws.Range("AB16").Value = 7
$ws.Range("AC16").Value = 0
$ws.Range("AD16").Value = "Rock"
$ws.Range("AE16").Value = "and"
$ws.Range("AF16").Value = "Desert"
$ws.Range("AG16").Value = "thousand"
$ws.Range("AH16").Value = 'km$^2$'

# Row 17 - Shrubland
$ws.Range("AB17").Value = 8
$ws.Range("AC17").Value = -20.365371
$ws.Range("AD17").Value = "Shrubland"
$ws.Range("AE17").Value = "thousand"
$ws.Range("AF17").Value = 'km$^2$'

# Row 18 - Tundra
$ws.Range("AB18").Value = 9
$ws.Range("AC18").Value = 0
$ws.Range("AD18").Value = "Tundra"
$ws.Range("AE18").Value = "thousand"
$ws.Range("AF18").Value = 'km$^2$'

# Row 19 - Urban
$ws.Range("AB19").Value = 10
$ws.Range("AC19").Value = 0
$ws.Range("AD19").Value = "Urban"
$ws.Range("AE19").Value = "thousand"
$ws.Range("AF19").Value = 'km$^2$'

# Row 22 - Land for Food summary = Animal Feed (AC9) + Crops (AC11)
$ws.Range("AB22").Value = "Land for Food"
$ws.Range("AC22").Formula = "=AC9+AC11"

# Scroll / selection to match the author's final view position
$ws.Range("AC25").Select()
$excel.ActiveWindow.ScrollColumn = 3
